$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from column J (rows 2-9) into column K so the new
# column picks up the same styles as its neighbour.
$ws.Range("J2:J9").Copy()
$ws.Range("K2:K9").PasteSpecial(-4122)

# Fill in the new column K values (year 2021 data).
$ws.Range("K3").Value = 2021
$ws.Range("K4").Value = 295
$ws.Range("K5").Value = 163
$ws.Range("K6").Value = 268
$ws.Range("K7").Value = 155
$ws.Range("K8").Value = 27
$ws.Range("K9").Value = 8

$ws.Range("L5").Select()
